$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176477789878845
$ws.Range("B1").Value = 2.511827945709229
$ws.Range("C1").Value = 9.450204849243164
$ws.Range("D1").Value = 2.100223541259766
$ws.Range("E1").Value = 1.224712371826172
